$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value  = -21.853
$ws.Range("A18").Value = -22.095
$ws.Range("A20").Value = -20.511
$ws.Range("A27").Value = -22.01
$ws.Range("A69").Value = -21.52
$ws.Range("A76").Value = -20.157
$ws.Range("A82").Value = -22.067
